# Rename the inline logo pictures in the document's headers/footers.
#
#   - BTec_Logo-Orange pictures (in both headers)     : image2.jpg -> image1.jpg
#   - PearsonLogo pictures      (in both footers)      : image1.png -> image2.png
#
# InlineShape objects don't expose a writable .Name property directly in
# this host, so each picture is briefly promoted to a floating Shape
# (ConvertToShape), renamed, then converted back to an inline shape
# (ConvertToInlineShape) so the drawing stays wp:inline, matching the
# original layout.

$d = $word.ActiveDocument
$sec = $d.Sections.First

function Rename-FirstInlinePicture($range, $newName) {
    $ishp = $range.InlineShapes.Item(1)
    $shp = $ishp.ConvertToShape()
    $shp.Name = $newName
    [void]$shp.ConvertToInlineShape()
}

# Headers (BTEC logo): image2.jpg -> image1.jpg
Rename-FirstInlinePicture $sec.Headers.Item(1).Range "image1.jpg"
Rename-FirstInlinePicture $sec.Headers.Item(2).Range "image1.jpg"

# Footers (Pearson logo): image1.png -> image2.png
Rename-FirstInlinePicture $sec.Footers.Item(1).Range "image2.png"
Rename-FirstInlinePicture $sec.Footers.Item(2).Range "image2.png"
